$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PMTestData")

# Row 54
$ws.Range("A54").Value = 'test_createDigitalExtension'
$ws.Range("B54").Value = 'number_initiate -number 80000000000000000000 -numbertype ex,80000000000000000000,MiVoice 4225 (DBC225),1B-2-20-00,FirstName,LastName,KSEXE:DIR=80000000000000000000;,number_end -number 80000000000000000000 -numbertype ex'
$ws.Range("B54").WrapText = $true
$ws.Range("C54").Value = 'Y'
$ws.Range("D54").Clear()
$ws.Rows.Item(54).RowHeight = 58

# Row 55
$ws.Range("A55").Value = 'test_edit_digital_extension_toChange_CAT'
$ws.Range("B55").Value = 'number_initiate -number 80000000000000000000 -numbertype ex,80000000000000000000,1B-2-20-00,KSEXE:DIR=80000000000000000000;,number_end -number 80000000000000000000 -numbertype ex'
$ws.Range("B55").WrapText = $true
$ws.Range("C55").Value = 'Y'
$ws.Range("D55").Clear()
$ws.Rows.Item(55).RowHeight = 58

# Row 56
$ws.Range("A56").Value = 'test_edit_digital_extension_toChange_First_Last_Names'
$ws.Range("B56").Value = 'number_initiate -number 80000000000000000000 -numbertype ex,80000000000000000000,MiVoice 4225 (DBC225),1B-2-20-00,FirstName,LastName,KSEXE:DIR=80000000000000000000;,number_end -number 80000000000000000000 -numbertype ex'
$ws.Range("B56").WrapText = $true
$ws.Range("C56").Value = 'Y'
$ws.Range("D56").Clear()
$ws.Rows.Item(56).RowHeight = 58

# Row 57
$ws.Range("A57").Value = 'test_edit_digital_extension_toChange_PhoneModel'
$ws.Range("B57").Value = 'number_initiate -number 80000000000000000000 -numbertype ex,80000000000000000000,1B-2-20-00,MiVoice 4222 (DBC222),KSEXE:DIR=80000000000000000000;,number_end -number 80000000000000000000 -numbertype ex'
$ws.Range("B57").WrapText = $true
$ws.Range("C57").Value = 'Y'
$ws.Range("D57").Clear()
$ws.Rows.Item(57).RowHeight = 58

# Row 58
$ws.Range("A58").Value = 'test_edit_digital_extension_toSet_AgentPosition'
$ws.Range("B58").Value = 'number_initiate -number 80000000000000000000 -numbertype ex,80000000000000000000,1B-2-20-00,KSEXE:DIR=80000000000000000000;,number_end -number 80000000000000000000 -numbertype ex'
$ws.Range("B58").WrapText = $true
$ws.Range("C58").Value = 'Y'
$ws.Range("D58").ClearContents()
$ws.Range("D58").WrapText = $true
$ws.Rows.Item(58).RowHeight = 58

# Row 59
$ws.Range("A59").Value = 'test_edit_digital_extension_toSet_HotLine'
$ws.Range("B59").Value = 'number_initiate -number 80000000000000000000 -numbertype ex,80000000000000000000,1B-2-20-00,KSEXE:DIR=80000000000000000000;,number_end -number 80000000000000000000 -numbertype ex'
$ws.Range("B59").WrapText = $true
$ws.Range("C59").Value = 'Y'
$ws.Range("D59").Value = 'number_initiate -number 90000000000000000000 -numbertype ex,extension -i -d 90000000000000000000 -l 1 --csp 0,ip_extension -i -d 90000000000000000000,90000000000000000000,ip_extension -e -d 90000000000000000000,extension -e -d 90000000000000000000,number_end -number 90000000000000000000 -numbertype ex'
$ws.Range("D59").WrapText = $true
$ws.Rows.Item(59).RowHeight = 72.5

# Row 60
$ws.Range("A60").Value = 'test_edit_digital_extension_toDelayed_HotLine'
$ws.Range("B60").Value = 'number_initiate -number 80000000000000000000 -numbertype ex,80000000000000000000,1B-2-20-00,KSEXE:DIR=80000000000000000000;,number_end -number 80000000000000000000 -numbertype ex'
$ws.Range("B60").WrapText = $true
$ws.Range("C60").Value = 'Y'
$ws.Range("D60").Value = 'number_initiate -number 90000000000000000000 -numbertype ex,extension -i -d 90000000000000000000 -l 1 --csp 0,ip_extension -i -d 90000000000000000000,90000000000000000000,ip_extension -e -d 90000000000000000000,extension -e -d 90000000000000000000,number_end -number 90000000000000000000 -numbertype ex'
$ws.Range("D60").WrapText = $true
$ws.Rows.Item(60).RowHeight = 72.5

# Row 61
$ws.Range("A61").Value = 'test_delete_digitalExtension'
$ws.Range("B61").Value = 'number_initiate -number 80000000000000000000 -numbertype ex,80000000000000000000,1B-2-20-00,KSEXE:DIR=80000000000000000000;,number_end -number 80000000000000000000 -numbertype ex'
$ws.Range("B61").WrapText = $true
$ws.Range("C61").Value = 'Y'
$ws.Range("D61").Clear()
$ws.Rows.Item(61).RowHeight = 58

# Row 62
$ws.Range("A62").Value = 'test_swap_digitalEquipmentPositions'
$ws.Range("B62").Value = 'number_initiate -number 80000000000000000000..80000000000000000001 -numbertype ex,80000000000000000000,1B-2-20-00,80000000000000000001,1B-2-20-01,KSEXE:DIR=80000000000000000000;,KSEXE:DIR=80000000000000000001;,number_end -number 80000000000000000000..80000000000000000001 -numbertype ex'
$ws.Range("B62").WrapText = $true
$ws.Range("C62").Value = 'Y'
$ws.Range("D62").Clear()
$ws.Rows.Item(62).RowHeight = 101.5

# Row 63
$ws.Range("A63").Value = 'test_create_digitalExt_using_Template'
$ws.Range("B63").Value = 'number_initiate -number 80000000000000000000 -numbertype ex,80000000000000000000,1B-2-20-00,DigitalTemplate,KSEXE:DIR=80000000000000000000;,number_end -number 80000000000000000000 -numbertype ex'
$ws.Range("B63").WrapText = $true
$ws.Range("C63").Value = 'Y'
$ws.Range("D63").Clear()
$ws.Rows.Item(63).RowHeight = 58

# Row 64
$ws.Range("A64").Value = 'test_createUser_with_Digital_Extension'
$ws.Range("B64").Value = 'TestUsr1,Mitel@123,Mitel@gmail.com,MitelFirst,MitelSecond, Business1,Business2,+917975935256,+918105855417'
$ws.Range("B64").WrapText = $true
$ws.Range("C64").Value = 'Y'
$ws.Range("D64").Value = 'number_initiate -number 80000000000000000000 -numbertype ex,80000000000000000000,MiVoice 4225 (DBC225),1B-2-20-00,FirstName,LastName,KSEXE:DIR=80000000000000000000;,number_end -number 80000000000000000000 -numbertype ex'
$ws.Range("D64").WrapText = $true
$ws.Rows.Item(64).RowHeight = 58

# Row 65
$ws.Range("A65").Value = 'test_createUser_with_existing_Digital_Extension'
$ws.Range("B65").Value = 'TestUsr1,Mitel@123,Mitel@gmail.com,MitelFirst,MitelSecond, Business1,Business2,+917975935256,+918105855417'
$ws.Range("B65").WrapText = $true
$ws.Range("C65").Value = 'Y'
$ws.Range("D65").Value = 'number_initiate -number 80000000000000000000 -numbertype ex,80000000000000000000,1B-2-20-00,KSEXE:DIR=80000000000000000000;,number_end -number 80000000000000000000 -numbertype ex'
$ws.Range("D65").WrapText = $true
$ws.Rows.Item(65).RowHeight = 43.5

# Row 66
$ws.Range("A66").Value = 'test_createUser_with_DigitalExtension_usingTemplate'
$ws.Range("B66").Value = 'TestUsr1,Mitel@123,Mitel@gmail.com,MitelFirst,MitelSecond, Business1,Business2,+917975935256,+918105855417'
$ws.Range("B66").WrapText = $true
$ws.Range("C66").Value = 'Y'
$ws.Range("D66").Value = 'number_initiate -number 80000000000000000000 -numbertype ex,80000000000000000000,1B-2-20-00,DigitalTemplate,KSEXE:DIR=80000000000000000000;,number_end -number 80000000000000000000 -numbertype ex'
$ws.Range("D66").WrapText = $true
$ws.Rows.Item(66).RowHeight = 58

# Row 67
$ws.Range("A67").Value = 'test_editUser_and_Assign_Existing_DigitalExtension'
$ws.Range("B67").Value = 'TestUsr1,Mitel@123,Mitel@gmail.com,MitelFirst,MitelSecond, Business1,Business2,+917975935256,+918105855417'
$ws.Range("B67").WrapText = $true
$ws.Range("C67").Value = 'Y'
$ws.Range("D67").Value = 'number_initiate -number 80000000000000000000 -numbertype ex,80000000000000000000,1B-2-20-00,KSEXE:DIR=80000000000000000000;,number_end -number 80000000000000000000 -numbertype ex'
$ws.Range("D67").WrapText = $true
$ws.Rows.Item(67).RowHeight = 43.5

# Row 68
$ws.Range("A68").Value = 'test_editUser_and_remove_Existing_DigitalExtension'
$ws.Range("B68").Value = 'TestUsr1,Mitel@123,Mitel@gmail.com,MitelFirst,MitelSecond, Business1,Business2,+917975935256,+918105855417'
$ws.Range("B68").WrapText = $true
$ws.Range("C68").Value = 'Y'
$ws.Range("D68").Value = 'number_initiate -number 80000000000000000000 -numbertype ex,80000000000000000000,MiVoice 4225 (DBC225),1B-2-20-00,KSEXE:DIR=80000000000000000000;,number_end -number 80000000000000000000 -numbertype ex'
$ws.Range("D68").WrapText = $true
$ws.Rows.Item(68).RowHeight = 58

# Row 69
$ws.Range("A69").Value = 'test_editUser_to_create_DigitalExtension'
$ws.Range("B69").Value = 'TestUsr1,Mitel@123,Mitel@gmail.com,MitelFirst,MitelSecond, Business1,Business2,+917975935256,+918105855417'
$ws.Range("B69").WrapText = $true
$ws.Range("C69").Value = 'Y'
$ws.Range("D69").Value = 'number_initiate -number 80000000000000000000 -numbertype ex,80000000000000000000,MiVoice 4225 (DBC225),1B-2-20-00,KSEXE:DIR=80000000000000000000;,number_end -number 80000000000000000000 -numbertype ex'
$ws.Range("D69").WrapText = $true
$ws.Rows.Item(69).RowHeight = 58

# Row 70
$ws.Range("A70").Value = 'test_createAnalogExtension'
$ws.Range("B70").Value = 'number_initiate -number 80000000000000000000 -numbertype ex,80000000000000000000,1B-2-10-00,FirstName,LastName,EXTEE:DIR=80000000000000000000;,number_end -number 80000000000000000000 -numbertype ex'
$ws.Range("B70").WrapText = $true
$ws.Range("C70").Value = 'Y'
$ws.Range("D70").Clear()
$ws.Rows.Item(70).RowHeight = 58

# Row 71
$ws.Range("A71").Value = 'test_edit_analog_extension_toChange_CAT'
$ws.Range("B71").Value = 'number_initiate -number 80000000000000000000 -numbertype ex,80000000000000000000,1B-2-10-00,EXTEE:DIR=80000000000000000000;,number_end -number 80000000000000000000 -numbertype ex'
$ws.Range("B71").WrapText = $true
$ws.Range("C71").Value = 'Y'
$ws.Range("D71").Clear()
$ws.Rows.Item(71).RowHeight = 58

# Row 72
$ws.Range("A72").Value = 'test_edit_analog_extension_toChange_First_Last_Names'
$ws.Range("B72").Value = 'number_initiate -number 80000000000000000000 -numbertype ex,80000000000000000000,1B-2-10-00,FirstName,LastName,EXTEE:DIR=80000000000000000000;,number_end -number 80000000000000000000 -numbertype ex'
$ws.Range("B72").WrapText = $true
$ws.Range("C72").Value = 'Y'
$ws.Range("D72").Clear()
$ws.Rows.Item(72).RowHeight = 58

# Row 73
$ws.Range("A73").Value = 'test_edit_analog_extension_equipmentPosition'
$ws.Range("B73").Value = 'number_initiate -number 80000000000000000000 -numbertype ex,80000000000000000000,1B-2-10-00,1B-2-10-01,EXTEE:DIR=80000000000000000000;,number_end -number 80000000000000000000 -numbertype ex'
$ws.Range("B73").WrapText = $true
$ws.Range("C73").Value = 'Y'
$ws.Range("D73").Clear()
$ws.Rows.Item(73).RowHeight = 58

# Row 74
$ws.Range("A74").Value = 'test_edit_analog_extension_toSet_HotLine'
$ws.Range("B74").Value = 'number_initiate -number 80000000000000000000 -numbertype ex,80000000000000000000,1B-2-10-00,EXTEE:DIR=80000000000000000000;,number_end -number 80000000000000000000 -numbertype ex'
$ws.Range("B74").WrapText = $true
$ws.Range("C74").Value = 'Y'
$ws.Range("D74").Value = 'number_initiate -number 90000000000000000000 -numbertype ex,extension -i -d 90000000000000000000 -l 1 --csp 0,ip_extension -i -d 90000000000000000000,90000000000000000000,ip_extension -e -d 90000000000000000000,extension -e -d 90000000000000000000,number_end -number 90000000000000000000 -numbertype ex'
$ws.Range("D74").WrapText = $true
$ws.Rows.Item(74).RowHeight = 72.5

# Row 75
$ws.Range("A75").Value = 'test_edit_analog_extension_toDelayed_HotLine'
$ws.Range("B75").Value = 'number_initiate -number 80000000000000000000 -numbertype ex,80000000000000000000,1B-2-10-00,EXTEE:DIR=80000000000000000000;,number_end -number 80000000000000000000 -numbertype ex'
$ws.Range("B75").WrapText = $true
$ws.Range("C75").Value = 'Y'
$ws.Range("D75").Value = 'number_initiate -number 90000000000000000000 -numbertype ex,extension -i -d 90000000000000000000 -l 1 --csp 0,ip_extension -i -d 90000000000000000000,90000000000000000000,ip_extension -e -d 90000000000000000000,extension -e -d 90000000000000000000,number_end -number 90000000000000000000 -numbertype ex'
$ws.Range("D75").WrapText = $true
$ws.Rows.Item(75).RowHeight = 72.5

# Row 76
$ws.Range("A76").Value = 'test_delete_analogExtension'
$ws.Range("B76").Value = 'number_initiate -number 80000000000000000000 -numbertype ex,80000000000000000000,1B-2-10-00,EXTEE:DIR=80000000000000000000;,number_end -number 80000000000000000000 -numbertype ex'
$ws.Range("B76").WrapText = $true
$ws.Range("C76").Value = 'Y'
$ws.Range("D76").Clear()
$ws.Rows.Item(76).RowHeight = 58

# Row 77
$ws.Range("A77").Value = 'test_swap_analogEquipmentPositions'
$ws.Range("B77").Value = 'number_initiate -number 80000000000000000000..80000000000000000001 -numbertype ex,80000000000000000000,1B-2-10-00,80000000000000000001,1B-2-10-01,EXTEE:DIR=80000000000000000000;,EXTEE:DIR=80000000000000000001;,number_end -number 80000000000000000000..80000000000000000001 -numbertype ex'
$ws.Range("B77").WrapText = $true
$ws.Range("C77").Value = 'Y'
$ws.Range("D77").Clear()
$ws.Rows.Item(77).RowHeight = 101.5

# Row 78
$ws.Range("A78").Value = 'test_create_analogExt_using_Template'
$ws.Range("B78").Value = 'number_initiate -number 80000000000000000000 -numbertype ex,80000000000000000000,1B-2-10-00,AnalogTemplate,EXTEE:DIR=80000000000000000000;,number_end -number 80000000000000000000 -numbertype ex'
$ws.Range("B78").WrapText = $true
$ws.Range("C78").Value = 'Y'
$ws.Range("D78").Clear()
$ws.Rows.Item(78).RowHeight = 58

# Row 79
$ws.Range("A79").Value = 'test_createUser_with_analog_Extension'
$ws.Range("B79").Value = 'TestUsr1,Mitel@123,Mitel@gmail.com,MitelFirst,MitelSecond, Business1,Business2,+917975935256,+918105855417'
$ws.Range("B79").WrapText = $true
$ws.Range("C79").Value = 'Y'
$ws.Range("D79").Value = 'number_initiate -number 80000000000000000000 -numbertype ex,80000000000000000000,1B-2-10-00,FirstName,LastName,EXTEE:DIR=80000000000000000000;,number_end -number 80000000000000000000 -numbertype ex'
$ws.Range("D79").WrapText = $true
$ws.Rows.Item(79).RowHeight = 58

# Row 80
$ws.Range("A80").Value = 'test_createUser_with_existing_analog_Extension'
$ws.Range("B80").Value = 'TestUsr1,Mitel@123,Mitel@gmail.com,MitelFirst,MitelSecond, Business1,Business2,+917975935256,+918105855417'
$ws.Range("B80").WrapText = $true
$ws.Range("C80").Value = 'Y'
$ws.Range("D80").Value = 'number_initiate -number 80000000000000000000 -numbertype ex,80000000000000000000,1B-2-10-00,EXTEE:DIR=80000000000000000000;,number_end -number 80000000000000000000 -numbertype ex'
$ws.Range("D80").WrapText = $true
$ws.Rows.Item(80).RowHeight = 43.5

# Row 81
$ws.Range("A81").Value = 'test_createUser_with_AnalogExtension_usingTemplate'
$ws.Range("B81").Value = 'TestUsr1,Mitel@123,Mitel@gmail.com,MitelFirst,MitelSecond, Business1,Business2,+917975935256,+918105855417'
$ws.Range("B81").WrapText = $true
$ws.Range("C81").Value = 'Y'
$ws.Range("D81").Value = 'number_initiate -number 80000000000000000000 -numbertype ex,80000000000000000000,1B-2-10-00,AnalogTemplate,EXTEE:DIR=80000000000000000000;,number_end -number 80000000000000000000 -numbertype ex'
$ws.Range("D81").WrapText = $true
$ws.Rows.Item(81).RowHeight = 58

# Row 82
$ws.Range("A82").Value = 'test_editUser_and_Assign_Existing_AnalogExtension'
$ws.Range("B82").Value = 'TestUsr1,Mitel@123,Mitel@gmail.com,MitelFirst,MitelSecond, Business1,Business2,+917975935256,+918105855417'
$ws.Range("B82").WrapText = $true
$ws.Range("C82").Value = 'Y'
$ws.Range("D82").Value = 'number_initiate -number 80000000000000000000 -numbertype ex,80000000000000000000,1B-2-10-00,EXTEE:DIR=80000000000000000000;,number_end -number 80000000000000000000 -numbertype ex'
$ws.Range("D82").WrapText = $true
$ws.Rows.Item(82).RowHeight = 43.5

# Row 83
$ws.Range("A83").Value = 'test_editUser_and_remove_Existing_AnalogExtension'
$ws.Range("B83").Value = 'TestUsr1,Mitel@123,Mitel@gmail.com,MitelFirst,MitelSecond, Business1,Business2,+917975935256,+918105855417'
$ws.Range("B83").WrapText = $true
$ws.Range("C83").Value = 'Y'
$ws.Range("D83").Value = 'number_initiate -number 80000000000000000000 -numbertype ex,80000000000000000000,1B-2-10-00,EXTEE:DIR=80000000000000000000;,number_end -number 80000000000000000000 -numbertype ex'
$ws.Range("D83").WrapText = $true
$ws.Rows.Item(83).RowHeight = 43.5

# Row 84
$ws.Range("A84").Value = 'test_editUser_to_create_AnalogExtension'
$ws.Range("B84").Value = 'TestUsr1,Mitel@123,Mitel@gmail.com,MitelFirst,MitelSecond, Business1,Business2,+917975935256,+918105855417'
$ws.Range("B84").WrapText = $true
$ws.Range("C84").Value = 'Y'
$ws.Range("D84").Value = 'number_initiate -number 80000000000000000000 -numbertype ex,80000000000000000000,1B-2-10-00,EXTEE:DIR=80000000000000000000;,number_end -number 80000000000000000000 -numbertype ex'
$ws.Range("D84").WrapText = $true
$ws.Rows.Item(84).RowHeight = 43.5
$ws.Activate()
$ws.Range("D81").Select()
$excel.ActiveWindow.ScrollRow = 79
$excel.ActiveWindow.ScrollColumn = 2
